# Apply "Taking latest code and Adding my changes" edits to the
# "Test Cases" sheet of the D suite workbook.
#
# Summary of the change:
#   - Runmode (col C) flips from "N" to "Y" for the test rows that were
#     previously skipped (rows 2-11).
#   - Results (col D) is reset to "SKIP" for the row that used to carry a
#     stale "PASS" (row 11, ProfileFollowingOthersTest) and for the two
#     rows that had no result yet (rows 12 and 13's prior state), while the
#     most-recently-run test (row 13, ProfileFollowerTest) now shows "PASS".
#   - The active selection left on the sheet moves to C12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Flip Runmode from N to Y for rows 2-11.
$ws.Range("C2:C11").Value = "Y"

# Update the Results column for rows 11-13.
$ws.Range("D11").Value = "SKIP"
$ws.Range("D12").Value = "SKIP"
$ws.Range("D13").Value = "PASS"

# Leave the selection where the author last clicked.
[void]$ws.Range("C12").Select()
